$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row right above the current row 65 ("44389" record),
# shifting every following row (65-185) down by one, so the sheet grows
# from 185 to 186 rows. Seed the new row by duplicating the row above it
# (row 64) so all the constant/label columns and the date-cell style come
# along for free, then overwrite the record-specific columns with the
# new observation's values.
$ws.Rows.Item(64).Copy()
$ws.Rows.Item(65).Insert()

$ws.Range("D65").Value = 44725
$ws.Range("J65").Value = 700
$ws.Range("K65").Value = 6000
$ws.Range("L65").Value = 6000
$ws.Range("M65").Value = 6000
$ws.Range("N65").Value = "`$/docena de matas"
$ws.Range("P65").Value = 1000
$ws.Range("Q65").Value = 6
